# Create the new "LeaveType" worksheet as the last tab in the workbook and
# populate it, matching the target diff:
#   - new sheet "LeaveType" appended after "Designation"
#   - becomes the active/selected tab (activeTab moves to it, tabSelected
#     moves off "Affiliate")
#   - two columns of leave data with custom column widths
#   - selection left on C23 on the new sheet

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end of
# the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$leaveType = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$leaveType.Name = "LeaveType"

# Populate header row.
$leaveType.Range("A1").Value = "Leave Name"
$leaveType.Range("B1").Value = "Leave Abbrevation"

# Populate data rows (abbreviation entered before the full name on row 2,
# matching the original authoring order so shared-string ids line up).
$leaveType.Range("B2").Value = "CL"
$leaveType.Range("A2").Value = "Casual Leave"

$leaveType.Range("A3").Value = "Sick Leave"
$leaveType.Range("B3").Value = "SL"

# Column widths matching the authored sheet.
$leaveType.Columns.Item(1).ColumnWidth = 24.66
$leaveType.Columns.Item(2).ColumnWidth = 25.16666

# Leave the on-sheet selection where the author left it.
$leaveType.Range("C23").Select()

Write-Output "LeaveType sheet created"
